$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "pc" row (row 23) entirely - shifts all subsequent rows up by one
$ws.Rows.Item(23).Delete()

# Append a new row of data at the bottom (new last row, 48)
$ws.Range("A48").Value = "zy_r"
$ws.Range("B48").Value = -66.85
$ws.Range("C48").Value = 12.19
$ws.Range("D48").Value = 50.46
